$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -0.2683900306044935
$ws.Range("D2").Value = 0.7908997231467618

# Row 3
$ws.Range("C3").Value = -0.81145292881394
$ws.Range("D3").Value = 0.4257982683053081

# Row 4
$ws.Range("C4").Value = -2.466225374689981
$ws.Range("D4").Value = 0.02192211205031702

# Row 5
$ws.Range("C5").Value = -3.943069974229491
$ws.Range("D5").Value = 0.0006928758670583512

# Row 6
$ws.Range("C6").Value = -0.4741569796156098
$ws.Range("D6").Value = 0.6400601932473062

# Row 7
$ws.Range("C7").Value = -2.204906921200918
$ws.Range("D7").Value = 0.0382152013885364

# Row 8
$ws.Range("C8").Value = -3.222872083801362
$ws.Range("D8").Value = 0.003915142416773154

# Row 9
$ws.Range("C9").Value = -1.603938226072133
$ws.Range("D9").Value = 0.1229886866583949
$ws.Range("G9").Value = "No"

# Row 10
$ws.Range("C10").Value = -2.503774378179249
$ws.Range("D10").Value = 0.02019990956805739

# Row 11
$ws.Range("C11").Value = -0.625074214271514
$ws.Range("D11").Value = 0.5383524450158101
